$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4881.3887
$ws.Range("I132").Value = 2565.3572
$ws.Range("J132").Value = 12987.5
$ws.Range("K132").Value = 7696.071599999999
$ws.Range("L132").Value = 38962.5
$ws.Range("M132").Value = -5166.071599999999
$ws.Range("N132").Value = -44022.5
$ws.Range("H138").Value = 1929.1464
$ws.Range("I138").Value = 1112.5641
$ws.Range("J138").Value = 2669.7673
$ws.Range("K138").Value = 3337.6923
$ws.Range("L138").Value = 8009.3019
$ws.Range("M138").Value = 1802.3077
$ws.Range("N138").Value = -18289.3019

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9584.138000000001
$ws.Range("I32").Value = 4677.04
$ws.Range("J32").Value = 40253.5
$ws.Range("K32").Value = 4677.04
$ws.Range("L32").Value = 40253.5
$ws.Range("M32").Value = -4390.04
$ws.Range("N32").Value = -40827.5
$ws.Range("H88").Value = 3089.1
$ws.Range("I88").Value = 3781.8333
$ws.Range("J88").Value = 2050
$ws.Range("K88").Value = 3781.8333
$ws.Range("L88").Value = 2050
$ws.Range("M88").Value = -3375.8333
$ws.Range("N88").Value = -2862
$ws.Range("H91").Value = 3089.1
$ws.Range("I91").Value = 3781.8333
$ws.Range("J91").Value = 2050
$ws.Range("K91").Value = 3781.8333
$ws.Range("L91").Value = 2050
$ws.Range("M91").Value = -2377.8333
$ws.Range("N91").Value = -4858
$ws.Range("H132").Value = 3127.681
$ws.Range("I132").Value = 2848.5757
$ws.Range("J132").Value = 3785.5715
$ws.Range("K132").Value = 8545.7271
$ws.Range("L132").Value = 11356.7145
$ws.Range("M132").Value = -6015.7271
$ws.Range("N132").Value = -16416.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7848.625
$ws.Range("I86").Value = 6607.8
$ws.Range("J86").Value = 9916.666999999999
$ws.Range("K86").Value = 6607.8
$ws.Range("L86").Value = 9916.666999999999
$ws.Range("M86").Value = -5484.8
$ws.Range("N86").Value = -12162.667
$ws.Range("H89").Value = 7848.625
$ws.Range("I89").Value = 6607.8
$ws.Range("J89").Value = 9916.666999999999
$ws.Range("K89").Value = 33039
$ws.Range("L89").Value = 49583.335
$ws.Range("M89").Value = -27423
$ws.Range("N89").Value = -60815.335
$ws.Range("H134").Value = 24911.156
$ws.Range("I134").Value = 30567.705
$ws.Range("J134").Value = 7427.273
$ws.Range("K134").Value = 91703.11500000001
$ws.Range("L134").Value = 22281.819
$ws.Range("M134").Value = -89168.11500000001
$ws.Range("N134").Value = -27351.819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2660.7534
$ws.Range("I31").Value = 1834.4419
$ws.Range("J31").Value = 3845.1333
$ws.Range("K31").Value = 1834.4419
$ws.Range("L31").Value = 3845.1333
$ws.Range("M31").Value = -1539.4419
$ws.Range("N31").Value = -4435.1333
$ws.Range("H34").Value = 2660.7534
$ws.Range("I34").Value = 1834.4419
$ws.Range("J34").Value = 3845.1333
$ws.Range("K34").Value = 1834.4419
$ws.Range("L34").Value = 3845.1333
$ws.Range("M34").Value = -1632.4419
$ws.Range("N34").Value = -4249.1333
$ws.Range("H58").Value = 2881.2307
$ws.Range("I58").Value = 2264.7273
$ws.Range("J58").Value = 3333.3333
$ws.Range("K58").Value = 2264.7273
$ws.Range("L58").Value = 3333.3333
$ws.Range("M58").Value = -2061.7273
$ws.Range("N58").Value = -3739.3333
$ws.Range("H99").Value = 86251
$ws.Range("I99").Value = 127989
$ws.Range("J99").Value = 2775
$ws.Range("K99").Value = 127989
$ws.Range("L99").Value = 2775
$ws.Range("M99").Value = -126491
$ws.Range("N99").Value = -5771
$ws.Range("H126").Value = 86251
$ws.Range("I126").Value = 127989
$ws.Range("J126").Value = 2775
$ws.Range("K126").Value = 383967
$ws.Range("L126").Value = 8325
$ws.Range("M126").Value = -381497
$ws.Range("N126").Value = -13265
$ws.Range("H132").Value = 2246.3264
$ws.Range("I132").Value = 1308.6
$ws.Range("J132").Value = 3726.9473
$ws.Range("K132").Value = 3925.8
$ws.Range("L132").Value = 11180.8419
$ws.Range("M132").Value = -1395.8
$ws.Range("N132").Value = -16240.8419
$ws.Range("H134").Value = 2040.3448
$ws.Range("I134").Value = 1411.3334
$ws.Range("K134").Value = 4234.0002
$ws.Range("M134").Value = -1699.0002
$ws.Range("H136").Value = 2881.2307
$ws.Range("I136").Value = 2264.7273
$ws.Range("J136").Value = 3333.3333
$ws.Range("K136").Value = 6794.1819
$ws.Range("L136").Value = 9999.999899999999
$ws.Range("M136").Value = -4244.1819
$ws.Range("N136").Value = -15099.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3316.25
$ws.Range("I3").Value = 3135.5833
$ws.Range("J3").Value = 3858.25
$ws.Range("K3").Value = 9406.749899999999
$ws.Range("L3").Value = 11574.75
$ws.Range("M3").Value = -9294.749899999999
$ws.Range("N3").Value = -11798.75
$ws.Range("H69").Value = 3131.55
$ws.Range("I69").Value = 777
$ws.Range("K69").Value = 2331
$ws.Range("M69").Value = -1520
$ws.Range("H72").Value = 3131.55
$ws.Range("I72").Value = 777
$ws.Range("K72").Value = 6993
$ws.Range("M72").Value = -2937
$ws.Range("H113").Value = 651.0769
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 651.0769
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1953.2307
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6293.2307
$ws.Range("H132").Value = 3480.925
$ws.Range("I132").Value = 1259.5333
$ws.Range("J132").Value = 4813.76
$ws.Range("K132").Value = 11335.7997
$ws.Range("L132").Value = 43323.84
$ws.Range("M132").Value = -8805.7997
$ws.Range("N132").Value = -48383.84

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3669.2888
$ws.Range("I132").Value = 3648.4827
$ws.Range("J132").Value = 3707
$ws.Range("K132").Value = 10945.4481
$ws.Range("L132").Value = 11121
$ws.Range("M132").Value = -8415.4481
$ws.Range("N132").Value = -16181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2003.9048
$ws.Range("I7").Value = 1882.3572
$ws.Range("J7").Value = 2247
$ws.Range("K7").Value = 1882.3572
$ws.Range("L7").Value = 2247
$ws.Range("M7").Value = -1770.3572
$ws.Range("N7").Value = -2471
$ws.Range("H46").Value = 593.8421
$ws.Range("I46").Value = 605.4
$ws.Range("J46").Value = 550.5
$ws.Range("K46").Value = 605.4
$ws.Range("L46").Value = 550.5
$ws.Range("M46").Value = -417.4
$ws.Range("N46").Value = -926.5
$ws.Range("H55").Value = 311.6
$ws.Range("I55").Value = 344.7857
$ws.Range("J55").Value = 269.36365
$ws.Range("K55").Value = 344.7857
$ws.Range("L55").Value = 269.36365
$ws.Range("M55").Value = -171.7857
$ws.Range("N55").Value = -615.36365
$ws.Range("H126").Value = 2003.9048
$ws.Range("I126").Value = 1882.3572
$ws.Range("J126").Value = 2247
$ws.Range("K126").Value = 5647.071599999999
$ws.Range("L126").Value = 6741
$ws.Range("M126").Value = -3177.071599999999
$ws.Range("N126").Value = -11681
$ws.Range("H132").Value = 7686.816
$ws.Range("I132").Value = 2563.9443
$ws.Range("J132").Value = 12297.4
$ws.Range("K132").Value = 7691.8329
$ws.Range("L132").Value = 36892.2
$ws.Range("M132").Value = -5161.8329
$ws.Range("N132").Value = -41952.2
$ws.Range("H136").Value = 4768.385
$ws.Range("I136").Value = 2410.1155
$ws.Range("J136").Value = 9484.923000000001
$ws.Range("K136").Value = 7230.3465
$ws.Range("L136").Value = 28454.769
$ws.Range("M136").Value = -4680.3465
$ws.Range("N136").Value = -33554.769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1276.8823
$ws.Range("I126").Value = 1113.8
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 3341.4
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -871.3999999999996
$ws.Range("N126").Value = -12440
$ws.Range("H132").Value = 1864.2167
$ws.Range("I132").Value = 1003.5952
$ws.Range("J132").Value = 3872.3333
$ws.Range("K132").Value = 3010.7856
$ws.Range("L132").Value = 11616.9999
$ws.Range("M132").Value = -480.7856000000002
$ws.Range("N132").Value = -16676.9999
$ws.Range("H136").Value = 25002370
$ws.Range("I136").Value = 71431030
$ws.Range("J136").Value = 2321.5386
$ws.Range("K136").Value = 214293090
$ws.Range("L136").Value = 6964.6158
$ws.Range("M136").Value = -214290540
$ws.Range("N136").Value = -12064.6158
